$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 27.02.2022 07:00"

# Update D8 (Delta Cena) from text "+1.6" to numeric value 1.6
$ws.Range("D8").Value = 1.6

# Update E8 (Old Datum) from text date to numeric date serial, formatted as date
$ws.Range("E8").Value = 44619.28143518518
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
